# Refresh cryptocurrency price/volume(1h) figures (GitHub Actions scrape update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'35.367.36"
$ws.Cells.Item(2, 5).Value = "'  +0.47%  "

$ws.Cells.Item(3, 4).Value = "'1.906.38"
$ws.Cells.Item(3, 5).Value = "'  +2.46%  "

$ws.Cells.Item(4, 5).Value = "'  -0.48%  "

$ws.Cells.Item(5, 4).Value = "'246.33"
$ws.Cells.Item(5, 5).Value = "'  +3.05%  "

$ws.Cells.Item(6, 4).Value = "'0.662"
$ws.Cells.Item(6, 5).Value = "'  +6.27%  "

$ws.Cells.Item(7, 5).Value = "'  -0.49%  "

$ws.Cells.Item(8, 4).Value = "'41.52"
$ws.Cells.Item(8, 5).Value = "'  -2.00%  "

$ws.Cells.Item(9, 4).Value = "'0.348"
$ws.Cells.Item(9, 5).Value = "'  +5.67%  "

$ws.Cells.Item(10, 4).Value = "'53.05"
$ws.Cells.Item(10, 5).Value = "'  +13.01%  "

$ws.Cells.Item(11, 4).Value = "'0.0721"
$ws.Cells.Item(11, 5).Value = "'  +3.98%  "

$ws.Cells.Item(12, 4).Value = "'0.0992"

$ws.Cells.Item(13, 4).Value = "'2.182.76"
$ws.Cells.Item(13, 5).Value = "'  +2.54%  "

$ws.Cells.Item(14, 5).Value = "'  +5.15%  "

$ws.Cells.Item(15, 5).Value = "'  +3.29%  "

$ws.Cells.Item(16, 4).Value = "'1.914.86"
$ws.Cells.Item(16, 5).Value = "'  +2.88%  "

$ws.Cells.Item(17, 4).Value = "'4.86"
$ws.Cells.Item(17, 5).Value = "'  +3.09%  "

$ws.Cells.Item(18, 4).Value = "'35.368.73"
$ws.Cells.Item(18, 5).Value = "'  +0.53%  "

$ws.Cells.Item(19, 4).Value = "'72.21"
$ws.Cells.Item(19, 5).Value = "'  +3.28%  "

$ws.Cells.Item(20, 4).Value = "'0.0₃0822"
$ws.Cells.Item(20, 5).Value = "'  +3.50%  "

$ws.Cells.Item(21, 4).Value = "'240.45"
$ws.Cells.Item(21, 5).Value = "'  -0.44%  "

$ws.Cells.Item(22, 4).Value = "'12.49"
$ws.Cells.Item(22, 5).Value = "'  +2.11%  "

$ws.Cells.Item(23, 4).Value = "'4.84"
$ws.Cells.Item(23, 5).Value = "'  +2.12%  "

$ws.Cells.Item(24, 5).Value = "'  -0.51%  "

$ws.Cells.Item(25, 4).Value = "'2.29"
$ws.Cells.Item(25, 5).Value = "'  +1.10%  "

$ws.Cells.Item(26, 5).Value = "'  +24.37%  "

$ws.Cells.Item(27, 4).Value = "'170.41"
$ws.Cells.Item(27, 5).Value = "'  +0.76%  "

$ws.Cells.Item(28, 5).Value = "'  +4.66%  "

$ws.Cells.Item(29, 4).Value = "'18.40"
$ws.Cells.Item(29, 5).Value = "'  +3.93%  "

$ws.Cells.Item(30, 5).Value = "'  +2.58%  "

$ws.Cells.Item(31, 5).Value = "'  +3.26%  "

$ws.Cells.Item(32, 5).Value = "'  +0.50%  "

$ws.Cells.Item(33, 5).Value = "'  +0.44%  "

$ws.Cells.Item(34, 4).Value = "'0.932"
$ws.Cells.Item(34, 5).Value = "'  +14.51%  "

$ws.Cells.Item(35, 4).Value = "'4.10"
$ws.Cells.Item(35, 5).Value = "'  +1.97%  "

$ws.Cells.Item(36, 4).Value = "'1.74"
$ws.Cells.Item(36, 5).Value = "'  -4.65%  "

$ws.Cells.Item(37, 4).Value = "'2.04"
$ws.Cells.Item(37, 5).Value = "'  -0.90%  "

$ws.Cells.Item(38, 5).Value = "'  +1.98%  "

$ws.Cells.Item(39, 5).Value = "'  +0.82%  "

$ws.Cells.Item(40, 4).Value = "'0.0208"
$ws.Cells.Item(40, 5).Value = "'  +3.02%  "

$ws.Cells.Item(41, 4).Value = "'16.37"
$ws.Cells.Item(41, 5).Value = "'  +8.40%  "

$ws.Cells.Item(42, 4).Value = "'0.0630"
$ws.Cells.Item(42, 5).Value = "'  +5.19%  "

$ws.Cells.Item(43, 4).Value = "'89.98"
$ws.Cells.Item(43, 5).Value = "'  -0.26%  "

$ws.Cells.Item(44, 4).Value = "'1.339.62"
$ws.Cells.Item(44, 5).Value = "'  -0.62%  "

$ws.Cells.Item(45, 5).Value = "'  +3.27%  "

$ws.Cells.Item(46, 4).Value = "'48.26"
$ws.Cells.Item(46, 5).Value = "'  +39.48%  "

$ws.Cells.Item(47, 2).Value = "'HuobiToken"
$ws.Cells.Item(47, 3).Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(47, 4).Value = "'2.40"
$ws.Cells.Item(47, 5).Value = "'  -0.85%  "

$ws.Cells.Item(48, 2).Value = "'MXToken"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(48, 4).Value = "'2.78"
$ws.Cells.Item(48, 5).Value = "'  +1.55%  "

$ws.Cells.Item(49, 4).Value = "'6.55"
$ws.Cells.Item(49, 5).Value = "'  -0.30%  "

$ws.Cells.Item(50, 4).Value = "'11.82"
$ws.Cells.Item(50, 5).Value = "'  -4.33%  "

$ws.Cells.Item(51, 4).Value = "'2.092.81"
$ws.Cells.Item(51, 5).Value = "'  +2.42%  "
